$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Message to the tutor"
$ws.Range("A2").Value = "Deji O"
$ws.Range("B2").Value = "I don’t like your teaching!!! I love it!!!"

$ws.Range("A3").Select()
